$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new article ("Central Mali: Top UN genocide prevention official sounds alarm
# over recent ethnically-targeted killings") was folded into the time-bucket
# analysis ahead of the existing "Mali arrests five suspects..." article, which
# pushes the two rows to swap places in the sheet (row 2 <-> row 3), while rows
# 4 and 5 stay exactly as they were.
#
# Only the displayed text of the uri column (E2/E3) is swapped here - the
# hyperlink targets that were already attached to E2/E3 are left untouched, so
# they keep pointing at their original articles (matching the source data
# pipeline's behaviour of re-writing cell text independently from the
# hyperlink relationships).

# Row 2 becomes the "Central Mali ..." article
$ws.Range("A2").Value = "Central Mali: Top UN genocide prevention official sounds alarm over recent ethnically-targeted killings"
$ws.Range("B2").Value = "2019-03-27T00:00:00UTC"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "day_2_to_30"
$ws.Range("E2").Value = "https://news.un.org/en/story/2019/03/1035661"

# Row 3 becomes the "Mali arrests five suspects ..." article
$ws.Range("A3").Value = "Mali arrests five suspects in killing of 157 villagers"
$ws.Range("B3").Value = "2019-03-30T10:20:36UTC"
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = "day_2_to_30"
$ws.Range("E3").Value = "https://www.africanews.com/2019/03/30/mali-arrests-five-suspects-in-killing-of-157-villagers/"
